$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "isMandatory" column (C) for boot_disk (row 7) and network_interface (row 8)
# from "yes" to "no"
$ws.Range("C7").Value = "no"
$ws.Range("C8").Value = "no"

# Update the visible scroll position / active selection on the sheet
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C6").Select()
